$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing existing rows 6-26 down to 7-27
$ws.Rows(6).Insert()

# Populate the new row 6 with the new data point
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44959
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 100112017
$ws.Range("G6").Value = "Ramas de apio"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = 5000
$ws.Range("N6").Value = "$/paquete"
$ws.Range("O6").Value = "Región de La Araucanía"
$ws.Range("P6").Value = 5000
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
